# Apply the edit described by the diff:
# - MTD sheet: update selection to C10 (single cell)
# - YTD sheet: update C10/D10 values, and update selection to C9 (single cell)

$wb = $excel.ActiveWorkbook

$wsMTD = $wb.Worksheets.Item("MTD")
$wsYTD = $wb.Worksheets.Item("YTD")

# Update data values on the YTD sheet (row 10: backlog)
$wsYTD.Range("C10").Value = 3800
$wsYTD.Range("D10").Value = 1200

# Update the saved cell selection on each sheet so the XML <selection> matches.
$wsMTD.Activate()
$wsMTD.Range("C10").Select() | Out-Null

$wsYTD.Activate()
$wsYTD.Range("C9").Select() | Out-Null
